# The "Folder Path" callout (Rectangle 16 / Straight Arrow Connector 17 /
# Rectangle 18) that was briefly added to the DATABASE STRUCTURE slide is
# removed again -- the path feature now auto-reads the absolute path instead
# of needing this manual annotation.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)

$s.Shapes.Item("Rectangle 16").Delete()
$s.Shapes.Item("Straight Arrow Connector 17").Delete()
$s.Shapes.Item("Rectangle 18").Delete()
